$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$base = "https://itos-humanitarian.s3.amazonaws.com"
$baseSlash = "https://itos-humanitarian.s3.amazonaws.com/"

# ---------------------------------------------------------------------------
# 1. New data rows 87-90: Bolivia (BOL) Admin0-Admin3 deployment entries
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 87; Base = "COD_BOL_Admin0" },
    @{ Row = 88; Base = "COD_BOL_Admin1" },
    @{ Row = 89; Base = "COD_BOL_Admin2" },
    @{ Row = 90; Base = "COD_BOL_Admin3" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $fileBase = $r.Base

    $ws.Cells.Item($row, 1).Value = 43980          # A: Date posted
    $ws.Cells.Item($row, 2).Value = $base           # B: Base URL for AWS
    $ws.Cells.Item($row, 3).Value = "BOL"           # C: country code
    $ws.Cells.Item($row, 4).Value = $fileBase       # D: base file name

    $urlGeo = "$baseSlash" + "BOL/" + $fileBase + ".geojson"
    $urlTopo = "$baseSlash" + "BOL/" + $fileBase + ".topojson"
    $urlKml = "$baseSlash" + "BOL/" + $fileBase + ".kml"
    $urlCsv = "$baseSlash" + "BOL/" + $fileBase + ".csv"

    $ws.Cells.Item($row, 5).Value = $urlGeo
    $ws.Cells.Item($row, 6).Value = $urlTopo
    $ws.Cells.Item($row, 7).Value = $urlKml
    $ws.Cells.Item($row, 8).Value = $urlCsv

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $baseSlash, "", "", $urlGeo)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $baseSlash, "", "", $urlTopo)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $baseSlash, "", "", $urlKml)

    if ($row -eq 89 -or $row -eq 90) {
        # H89/H90 keep their literal value as display text (no explicit display passed)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 8), $baseSlash, "", "", "")
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 8), $baseSlash, "", "", $urlCsv)
    }
}

# ---------------------------------------------------------------------------
# 2. Cell formatting: wrap text + vertical-center for the new URL cells
#    (built on a scratch cell first so the style transform is applied in a
#    single diff, matching the two new cellXfs entries added upstream)
# ---------------------------------------------------------------------------

$scratchPlain = $ws.Range("Z1")
$scratchPlain.WrapText = $true
$scratchPlain.VerticalAlignment = -4108
$scratchPlain.Copy()
$ws.Range("G91:H91").PasteSpecial(-4122)
$scratchPlain.Clear()

$scratchLink = $ws.Range("Z1")
$scratchLink.Value = "x"
$ws.Hyperlinks.Add($scratchLink, $baseSlash, "", "", "")
$scratchLink.WrapText = $true
$scratchLink.VerticalAlignment = -4108
$scratchLink.Copy()
$ws.Range("E87:H90").PasteSpecial(-4122)
$scratchLink.Clear()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. View state: land the sheet roughly where the author left it
# ---------------------------------------------------------------------------

$excel.Goto($ws.Range("C91"), $true)

Write-Output "done"
